$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 67-68; this pushes the existing rows 67-181
# down to 69-183 (matching the target diff, where old row 67 becomes
# new row 69, ... old row 181 becomes new row 183).
$ws.Rows("67:68").Insert()

# Populate the two freshly inserted rows with the new record data.
# Columns A,B,C,E,F,G,H,I,J,K,Q,R,T are constant across this whole
# data block, so just replicate them; only D,L,M,N,O,P,S differ.

# Row 67: "Especial" quality record
$ws.Range("A67").Value = 3
$ws.Range("B67").Value = "Femacal de La Calera"
$ws.Range("C67").Value = "Coquimbo"
$ws.Range("D67").Value = 44533
$ws.Range("E67").Value = 5
$ws.Range("F67").Value = "Fruta"
$ws.Range("G67").Value = 100101
$ws.Range("H67").Value = "Berries"
$ws.Range("I67").Value = 100112025
$ws.Range("J67").Value = "Frutilla"
$ws.Range("K67").Value = "Sin especificar"
$ws.Range("L67").Value = "Especial"
$ws.Range("M67").Value = 139
$ws.Range("N67").Value = 5500
$ws.Range("O67").Value = 6000
$ws.Range("P67").Value = 5806
$ws.Range("Q67").Value = "`$/bandeja 7 kilos"
$ws.Range("R67").Value = "Provincia de Melipilla"
$ws.Range("S67").Value = 829
$ws.Range("T67").Value = 7

# Row 68: "Segunda" quality record
$ws.Range("A68").Value = 3
$ws.Range("B68").Value = "Femacal de La Calera"
$ws.Range("C68").Value = "Coquimbo"
$ws.Range("D68").Value = 44533
$ws.Range("E68").Value = 5
$ws.Range("F68").Value = "Fruta"
$ws.Range("G68").Value = 100101
$ws.Range("H68").Value = "Berries"
$ws.Range("I68").Value = 100112025
$ws.Range("J68").Value = "Frutilla"
$ws.Range("K68").Value = "Sin especificar"
$ws.Range("L68").Value = "Segunda"
$ws.Range("M68").Value = 78
$ws.Range("N68").Value = 4000
$ws.Range("O68").Value = 4000
$ws.Range("P68").Value = 4000
$ws.Range("Q68").Value = "`$/bandeja 7 kilos"
$ws.Range("R68").Value = "Provincia de Melipilla"
$ws.Range("S68").Value = 571
$ws.Range("T68").Value = 7
